$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.604.72"
$ws.Range("E2").Value = "  +2.03%  "
$ws.Range("D3").Value = "1.844.98"
$ws.Range("E3").Value = "  +4.01%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.65%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.555"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.57%  "
$ws.Range("E7").Value = "  -0.70%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.78"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.84%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.294"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.28%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0715"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +9.37%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0932"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.37%  "
$ws.Range("D12").Value = "2.109.22"
$ws.Range("E12").Value = "  +3.84%  "
$ws.Range("E13").Value = "  +1.20%  "
$ws.Range("D14").Value = "1.843.56"
$ws.Range("E14").Value = "  +3.90%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.653"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.75%  "
$ws.Range("D16").Value = "34.651.32"
$ws.Range("E16").Value = "  +2.08%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.36"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.11%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "70.14"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "255.08"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.76%  "
$ws.Range("D20").Value = "0.0₃0807"
$ws.Range("E20").Value = "  +9.84%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.32"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +9.87%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.997"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.62%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.32"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.37%  "
$ws.Range("E24").Value = "  +1.57%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "161.68"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.54%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "17.04"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.25"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.44%  "
$ws.Range("E28").Value = "  +2.04%  "
$ws.Range("E29").Value = "  -0.73%  "
$ws.Range("E30").Value = "  +3.78%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.84"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.68%  "
$ws.Range("E32").Value = "  +1.56%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "504.40"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +864.90%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.64"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.60%  "
$ws.Range("E35").Value = "  +7.05%  "
$ws.Range("D36").Value = "1.446.73"
$ws.Range("E36").Value = "  +0.23%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.659"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.70%  "
$ws.Range("E38").Value = "  +1.63%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0194"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.40%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.984"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +11.59%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "83.19"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.26%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.82"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.83%  "
$ws.Range("E43").Value = "  +0.33%  "
$ws.Range("E44").Value = "  +6.49%  "
$ws.Range("E45").Value = "  +6.27%  "
$ws.Range("D46").Value = "1.999.27"
$ws.Range("E46").Value = "  +3.48%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.56"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +7.34%  "
$ws.Range("E48").Value = "  -0.12%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0492"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.17%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "106.79"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +10.07%  "
